$wb = $excel.ActiveWorkbook

# --- Sheet4 (tab position 3, "Sheet4") updates: trucking detail additions ---
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Range("B19").Value = "0000 for auto gen"
$ws4.Range("F13").Value = "LCL -> add -> FCL :bug"

# --- Sheet6 (tab position 6, "Sheet6") updates: fix report bug for sales ---
$ws6 = $wb.Worksheets.Item("Sheet6")
$ws6.Range("A1").ClearContents()
$ws6.Range("A2").ClearContents()
$ws6.Range("C2").ClearContents()
$ws6.Range("A8").ClearContents()
$ws6.Range("E8").ClearContents()
$ws6.Range("H15").Value = "4. Hoi lai phieu de nghi thanh toan"
$ws6.Range("H16").Value = "5. Phieu thu chi ouref"
$ws6.Range("H17").Value = "refund cho 1 job???"
$ws6.Range("H18").Value = "6. Ke hoach van tai"
$ws6.Range("H19").Value = "7. Container template"

# --- Selection / active tab updates ---
$ws4.Activate()
$ws4.Range("H25").Select()

$ws6.Activate()
$ws6.Range("E8").Select()
